$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.807.87"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.745.87"
$ws.Range("E3").Value = "  -4.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.47"
$ws.Range("E5").Value = "  -8.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -5.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.96"
$ws.Range("E8").Value = "  -6.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2652"
$ws.Range("E9").Value = "  -11.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06167"
$ws.Range("E10").Value = "  -10.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.743.60"
$ws.Range("E11").Value = "  -5.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06929"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.43"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6004"
$ws.Range("E14").Value = "  -18.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.502"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.01"
$ws.Range("E16").Value = "  -13.94%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.810.06"
$ws.Range("E19").Value = "  -2.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006883"
$ws.Range("E20").Value = "  -12.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.966.24"
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.043"
$ws.Range("E23").Value = "  -11.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.212"
$ws.Range("E24").Value = "  -12.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.145"
$ws.Range("E25").Value = "  -11.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.79"
$ws.Range("E26").Value = "  -3.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.519"
$ws.Range("E27").Value = "  -9.87%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.817"
$ws.Range("E28").Value = "  -16.77%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.01"
$ws.Range("E29").Value = "  -11.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.59"
$ws.Range("E30").Value = "  -6.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.775"
$ws.Range("E31").Value = "  -10.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08111"
$ws.Range("E32").Value = "  -8.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.473"
$ws.Range("E33").Value = "  -13.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04518"
$ws.Range("E34").Value = "  -5.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9987"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.649"
$ws.Range("E36").Value = "  -9.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9843"
$ws.Range("E37").Value = "  -12.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6081"
$ws.Range("E38").Value = "  -16.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.669"
$ws.Range("E39").Value = "  -13.71%  "
$ws.Range("E40").Value = "  -9.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.913"
$ws.Range("E41").Value = "  -15.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.93"
$ws.Range("E43").Value = "  -4.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3815"
$ws.Range("E44").Value = "  -19.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.089"
$ws.Range("E45").Value = "  -13.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7342"
$ws.Range("E46").Value = "  -18.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05355"
$ws.Range("E47").Value = "  -7.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1111"
$ws.Range("E48").Value = "  -10.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.24"
$ws.Range("E49").Value = "  -12.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.924"
$ws.Range("E50").Value = "  -19.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.63"
$ws.Range("E51").Value = "  -12.18%  "
